# Commit: "Include clustering before EN model"
#
# A new calibration run (electrode E50, dated 2022-02-07) was logged on the
# "Low DA" sheet as row 13. Downstream artifacts (the "Count" sheet's
# COUNTIF tallies in column E, and the bar chart built on top of it) are
# formulas, so adding the row and letting Excel recalculate reproduces the
# rest of the diff automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Low DA")

# New row 13: Date | Electrode | sixteen sweep-voltage sample points (C:Q)
$ws.Range("A13").Value = 44599      # 2022-02-07
$ws.Range("B13").Value = "E50"

$samplePoints = @(5, 25, 35, 45, 50, 55, 60, 95, 100, 110, 115, 135, 140, 145, 150)
for ($i = 0; $i -lt $samplePoints.Length; $i++) {
    $ws.Cells.Item(13, 3 + $i).Value = $samplePoints[$i]
}

# Recalculate so the "Count" sheet's COUNTIF formulas (and the chart that
# reads from it) pick up the new row.
$excel.CalculateFullRebuild()
